$d = $word.ActiveDocument

# Update the date line at the top of the document.
$d.Content.Find.Execute("2024-04-12 Friday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2024-04-13 Saturday", 2)

# Update the division problems in the table. Several values repeat
# elsewhere in the table, so each replacement is scoped to its specific
# cell (row, column) to avoid touching the wrong occurrence.
$t = $d.Tables.Item(1)

$replacements = @(
    @{ Row = 1;  Col = 1; Old = "97÷2=48, 1";  New = "87÷5=17, 2" },
    @{ Row = 1;  Col = 2; Old = "70÷4=17, 2";  New = "22÷3=7, 1" },
    @{ Row = 1;  Col = 3; Old = "67÷5=13, 2";  New = "25÷3=8, 1" },
    @{ Row = 1;  Col = 4; Old = "35÷8=4, 3";   New = "54÷4=13, 2" },
    @{ Row = 1;  Col = 5; Old = "88÷3=29, 1";  New = "30÷5=6, 0" },

    @{ Row = 5;  Col = 1; Old = "48÷9=5, 3";   New = "43÷6=7, 1" },
    @{ Row = 5;  Col = 2; Old = "10÷6=1, 4";   New = "77÷6=12, 5" },
    @{ Row = 5;  Col = 3; Old = "27÷4=6, 3";   New = "51÷2=25, 1" },
    @{ Row = 5;  Col = 4; Old = "36÷8=4, 4";   New = "43÷2=21, 1" },
    @{ Row = 5;  Col = 5; Old = "53÷3=17, 2";  New = "60÷4=15, 0" },

    @{ Row = 9;  Col = 1; Old = "57÷8=7, 1";   New = "52÷3=17, 1" },
    @{ Row = 9;  Col = 2; Old = "80÷2=40, 0";  New = "18÷7=2, 4" },
    @{ Row = 9;  Col = 3; Old = "19÷9=2, 1";   New = "62÷8=7, 6" },
    @{ Row = 9;  Col = 4; Old = "35÷5=7, 0";   New = "55÷2=27, 1" },
    @{ Row = 9;  Col = 5; Old = "96÷8=12, 0";  New = "29÷8=3, 5" },

    @{ Row = 13; Col = 1; Old = "88÷3=29, 1";  New = "38÷7=5, 3" },
    @{ Row = 13; Col = 2; Old = "45÷7=6, 3";   New = "91÷6=15, 1" },
    @{ Row = 13; Col = 3; Old = "56÷5=11, 1";  New = "89÷2=44, 1" },
    @{ Row = 13; Col = 4; Old = "32÷9=3, 5";   New = "69÷3=23, 0" },
    @{ Row = 13; Col = 5; Old = "71÷4=17, 3";  New = "35÷5=7, 0" },

    @{ Row = 17; Col = 1; Old = "39÷8=4, 7";   New = "70÷8=8, 6" },
    @{ Row = 17; Col = 2; Old = "69÷2=34, 1";  New = "28÷4=7, 0" },
    @{ Row = 17; Col = 3; Old = "32÷9=3, 5";   New = "56÷8=7, 0" },
    @{ Row = 17; Col = 4; Old = "25÷7=3, 4";   New = "63÷9=7, 0" },
    @{ Row = 17; Col = 5; Old = "15÷9=1, 6";   New = "85÷6=14, 1" }
)

foreach ($item in $replacements) {
    # Find/Replace isn't reliably bounded to a single cell when several
    # cells share identical text, so set the cell Range's Text directly
    # instead -- this is scoped precisely to the targeted cell.
    $cellRange = $t.Cell($item.Row, $item.Col).Range
    $cellRange.Text = $item.New
}
